{"js": "// Update the cover-page version/date line:\n//   \"Version 11.00.00, 2015-03-16\"  ->  \"Version 11.01.00, 2015-03-23\"\nconst oldText = \"Version 11.00.00, 2015-03-16\";\nconst newText = \"Version 11.01.00, 2015-03-23\";\n\nconst results = context.document.body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Replace the matched range's text in place, preserving its formatting.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n} else {\n  // Fallback: the string may be split oddly; locate the paragraph that\n  // starts with \"Version \" and rebuild its text directly.\n  const paragraphs = context.document.body.paragraphs;\n  paragraphs.load(\"items/text\");\n  await context.sync();\n\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    const para = paragraphs.items[i];\n    if (para.text.indexOf(\"Version 11.00.00, 2015-03-16\") !== -1) {\n      const range = para.getRange();\n      range.insertText(\n        para.text.replace(\"Version 11.00.00, 2015-03-16\", newText),\n        Word.InsertLocation.replace\n      );\n      break;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the cover-page version/date line:\n#   \"Version 11.00.00, 2015-03-16\"  ->  \"Version 11.01.00, 2015-03-23\"\n\n$d = $word.ActiveDocument\n\n$oldText = \"Version 11.00.00, 2015-03-16\"\n$newText = \"Version 11.01.00, 2015-03-23\"\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\n    $oldText,   # FindText\n    $false,     # MatchCase\n    $false,     # MatchWholeWord\n    $false,     # MatchWildcards\n    $false,     # MatchSoundsLike\n    $false,     # MatchAllWordForms\n    $true,      # Forward\n    1,          # Wrap (wdFindContinue)\n    $false,     # Format\n    $newText,   # ReplaceWith\n    2           # Replace (wdReplaceAll)\n)\n\nif (-not $found) {\n    # Fallback: locate the paragraph that contains the old version string\n    # and replace its text directly.\n    foreach ($para in $d.Paragraphs) {\n        $pr = $para.Range\n        if ($pr.Text -like \"*$oldText*\") {\n            $pr.Text = $pr.Text.Replace($oldText, $newText)\n            break\n        }\n    }\n}\n"}
